$p = $ppt.ActivePresentation
Write-Host "Slide count:" $p.Slides.Count
